$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(4016, 4199, 4429, 4787, 5245, 5245, 5245, 5245, 5245, 5245, 5245, 5245, 5304, 5304)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
